$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. semantic_aspect_model_schema: drop the long "description" row, keep the
#    machine-readable header row (shifts from row 2 up to row 1), rename the
#    path segments to the indexed form, and re-style the row.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("semantic_aspect_model_schema")

$ws1.Rows.Item(1).Delete()

$ws1.Range("A1").Value2 = "catenaXId"
$ws1.Range("B1").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_latitude"
$ws1.Range("C1").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_longitude"
$ws1.Range("D1").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_altitude"
$ws1.Range("E1").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_geoDataTimestamp"
$ws1.Range("F1").Value2 = "sensorRuntimeData[0]_batteryLevel"
$ws1.Range("G1").Value2 = "sensorRuntimeData[0]_timestamp"
$ws1.Range("H1").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorType"
$ws1.Range("I1").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorValue"
$ws1.Range("J1").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorUnit"
$ws1.Range("K1").Value2 = "sensorRuntimeData[0]_transmissionMethod"

$hdrRange = $ws1.Range("A1:K1")
$hdrRange.WrapText = $False
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4108

# Column widths (values chosen so that Excel's internal pixel rounding lands
# on the target character widths: 10.8 / 39.6 / 36 / 46.8).
$ws1.Columns.Item(1).ColumnWidth = 9.95
$ws1.Columns.Item(6).ColumnWidth = 38.77
$ws1.Columns.Item(7).ColumnWidth = 35.15
$ws1.Columns.Item(11).ColumnWidth = 45.96

# ---------------------------------------------------------------------------
# 2. New "description" tab: Column Name / Description / Possible Values.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDesc = $wb.Worksheets.Add($null, $lastSheet)
$wsDesc.Name = "description"

$wsDesc.Range("A1").Value2 = "Column Name"
$wsDesc.Range("B1").Value2 = "Description"
$wsDesc.Range("C1").Value2 = "Possible Values"

$wsDesc.Range("A2").Value2 = "catenaXId"
$wsDesc.Range("B2").Value2 = "The fully anonymous Catena-X ID of the asset, valid for the Catena-X dataspace."

$wsDesc.Range("A3").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_latitude"
$wsDesc.Range("B3").Value2 = "The angle between zenith and a plane parallel to the equator."

$wsDesc.Range("A4").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_longitude"
$wsDesc.Range("B4").Value2 = "Geographic coordinate that specifies the east-west position of a point on the Earth's surface."

$wsDesc.Range("A5").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_altitude"
$wsDesc.Range("B5").Value2 = "Antenna Altitude above/below mean-sea-level (geoid)."

$wsDesc.Range("A6").Value2 = "sensorRuntimeData[0]_sensorGeoLocation_geoDataTimestamp"
$wsDesc.Range("B6").Value2 = "The timestamp of the latest sensor reading of the geo data."

$wsDesc.Range("A7").Value2 = "sensorRuntimeData[0]_batteryLevel"
$wsDesc.Range("B7").Value2 = "The battery level displays how much charge of the battery has been left."

$wsDesc.Range("A8").Value2 = "sensorRuntimeData[0]_timestamp"
$wsDesc.Range("B8").Value2 = "The timestamp of the latest sensor reading."

$wsDesc.Range("A9").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorType"
$wsDesc.Range("B9").Value2 = "Different types of sensors that are commonly used in various applications, measuring one of the physical properties like Temperature, Pressure,  Resistance, Shock, Conduction, Heat Transfer etc."

$wsDesc.Range("A10").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorValue"
$wsDesc.Range("B10").Value2 = "The measured value of the sensor type."

$wsDesc.Range("A11").Value2 = "sensorRuntimeData[0]_sensorData[0]_sensorUnit"
$wsDesc.Range("B11").Value2 = "Describes a Property containing a reference to one of the units in the Unit Catalog."

$wsDesc.Range("A12").Value2 = "sensorRuntimeData[0]_transmissionMethod"
$wsDesc.Range("B12").Value2 = "The method under which the sensing data is transmitted from the source to the remote node."

$wsDesc.Columns.Item(1).ColumnWidth = 49.13
$wsDesc.Columns.Item(2).ColumnWidth = 49.13
$wsDesc.Columns.Item(3).ColumnWidth = 16.13

# ---------------------------------------------------------------------------
# 3. New "metadata" tab (hidden): provenance of the auto-generated artifact.
# ---------------------------------------------------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMeta = $wb.Worksheets.Add($null, $lastSheet2)
$wsMeta.Name = "metadata"

$wsMeta.Range("A1").Value2 = "Key"
$wsMeta.Range("B1").Value2 = "Value"

$wsMeta.Range("A2").Value2 = "basedOnCommit"
$wsMeta.Range("B2").Value2 = "1b740b427f8155bb666b8855ec574b6765198bdb"

$wsMeta.Range("A3").Value2 = "commitHtmlUrl"
$wsMeta.Range("B3").Value2 = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/1b740b427f8155bb666b8855ec574b6765198bdb"

$wsMeta.Range("A4").Value2 = "commitDate"
$wsMeta.Range("B4").Value2 = "2023-11-21 09:49:45+00:00"

$wsMeta.Range("A5").Value2 = "commitMessage"
$wsMeta.Range("B5").Value2 = "Adding auto-generated artifacts for new models"

$wsMeta.Visible = $False

# ---------------------------------------------------------------------------
# 4. Restore the original active sheet / selection.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("A1").Select()
